$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$r = 307
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 307
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-arrows-dark-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'arrow navigation motion'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A307,",''title'':''",C307,"'',''group'':''",D307,"'',''area'':''",E307,"'',''keywords'':[''",SUBSTITUTE(F307," ","'',''"),"''],''description'':''",G307,"'',''publish'':",H307,"}"),"''","""")'

$r = 308
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 308
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-big-code-1-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'code'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A308,",''title'':''",C308,"'',''group'':''",D308,"'',''area'':''",E308,"'',''keywords'':[''",SUBSTITUTE(F308," ","'',''"),"''],''description'':''",G308,"'',''publish'':",H308,"}"),"''","""")'

$r = 309
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 309
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-big-code-2-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'code'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A309,",''title'':''",C309,"'',''group'':''",D309,"'',''area'':''",E309,"'',''keywords'':[''",SUBSTITUTE(F309," ","'',''"),"''],''description'':''",G309,"'',''publish'':",H309,"}"),"''","""")'

$r = 310
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 310
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-big-code-dark-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'code'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A310,",''title'':''",C310,"'',''group'':''",D310,"'',''area'':''",E310,"'',''keywords'':[''",SUBSTITUTE(F310," ","'',''"),"''],''description'':''",G310,"'',''publish'':",H310,"}"),"''","""")'

$r = 311
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 311
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-big-code-turquoise'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'code'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A311,",''title'':''",C311,"'',''group'':''",D311,"'',''area'':''",E311,"'',''keywords'':[''",SUBSTITUTE(F311," ","'',''"),"''],''description'':''",G311,"'',''publish'':",H311,"}"),"''","""")'

$r = 312
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 312
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-bugs-gray'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'bug debug'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A312,",''title'':''",C312,"'',''group'':''",D312,"'',''area'':''",E312,"'',''keywords'':[''",SUBSTITUTE(F312," ","'',''"),"''],''description'':''",G312,"'',''publish'':",H312,"}"),"''","""")'

$r = 313
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 313
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-bugs-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'bug debug'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A313,",''title'':''",C313,"'',''group'':''",D313,"'',''area'':''",E313,"'',''keywords'':[''",SUBSTITUTE(F313," ","'',''"),"''],''description'':''",G313,"'',''publish'':",H313,"}"),"''","""")'

$r = 314
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 314
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-bugs-turquoise'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'bug debug'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A314,",''title'':''",C314,"'',''group'':''",D314,"'',''area'':''",E314,"'',''keywords'':[''",SUBSTITUTE(F314," ","'',''"),"''],''description'':''",G314,"'',''publish'':",H314,"}"),"''","""")'

$r = 315
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 315
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-cloud-vectors-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'cloud'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A315,",''title'':''",C315,"'',''group'':''",D315,"'',''area'':''",E315,"'',''keywords'':[''",SUBSTITUTE(F315," ","'',''"),"''],''description'':''",G315,"'',''publish'':",H315,"}"),"''","""")'

$r = 316
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 316
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-code-1-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'code'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A316,",''title'':''",C316,"'',''group'':''",D316,"'',''area'':''",E316,"'',''keywords'':[''",SUBSTITUTE(F316," ","'',''"),"''],''description'':''",G316,"'',''publish'':",H316,"}"),"''","""")'

$r = 317
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 317
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-code-2-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'code'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A317,",''title'':''",C317,"'',''group'':''",D317,"'',''area'':''",E317,"'',''keywords'':[''",SUBSTITUTE(F317," ","'',''"),"''],''description'':''",G317,"'',''publish'':",H317,"}"),"''","""")'

$r = 318
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 318
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-cubes-magenta'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'cube microservice module component'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A318,",''title'':''",C318,"'',''group'':''",D318,"'',''area'':''",E318,"'',''keywords'':[''",SUBSTITUTE(F318," ","'',''"),"''],''description'':''",G318,"'',''publish'':",H318,"}"),"''","""")'

$r = 319
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 319
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-devices-1-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'device mobile phone computer laptop'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A319,",''title'':''",C319,"'',''group'':''",D319,"'',''area'':''",E319,"'',''keywords'':[''",SUBSTITUTE(F319," ","'',''"),"''],''description'':''",G319,"'',''publish'':",H319,"}"),"''","""")'

$r = 320
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 320
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-devices-1-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'device mobile phone computer laptop'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A320,",''title'':''",C320,"'',''group'':''",D320,"'',''area'':''",E320,"'',''keywords'':[''",SUBSTITUTE(F320," ","'',''"),"''],''description'':''",G320,"'',''publish'':",H320,"}"),"''","""")'

$r = 321
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 321
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-mix-1-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'icon'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A321,",''title'':''",C321,"'',''group'':''",D321,"'',''area'':''",E321,"'',''keywords'':[''",SUBSTITUTE(F321," ","'',''"),"''],''description'':''",G321,"'',''publish'':",H321,"}"),"''","""")'

$r = 322
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 322
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-mix-1-gray'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'icon'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A322,",''title'':''",C322,"'',''group'':''",D322,"'',''area'':''",E322,"'',''keywords'':[''",SUBSTITUTE(F322," ","'',''"),"''],''description'':''",G322,"'',''publish'':",H322,"}"),"''","""")'

$r = 323
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 323
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-outline-icons-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'icon'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A323,",''title'':''",C323,"'',''group'':''",D323,"'',''area'':''",E323,"'',''keywords'':[''",SUBSTITUTE(F323," ","'',''"),"''],''description'':''",G323,"'',''publish'':",H323,"}"),"''","""")'

$r = 324
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 324
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-outline-icons-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'particle line'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A324,",''title'':''",C324,"'',''group'':''",D324,"'',''area'':''",E324,"'',''keywords'':[''",SUBSTITUTE(F324," ","'',''"),"''],''description'':''",G324,"'',''publish'':",H324,"}"),"''","""")'

$r = 325
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 325
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-particle-lines-blue'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'particle line'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A325,",''title'':''",C325,"'',''group'':''",D325,"'',''area'':''",E325,"'',''keywords'':[''",SUBSTITUTE(F325," ","'',''"),"''],''description'':''",G325,"'',''publish'':",H325,"}"),"''","""")'

$r = 326
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 326
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-people-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'people'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A326,",''title'':''",C326,"'',''group'':''",D326,"'',''area'':''",E326,"'',''keywords'':[''",SUBSTITUTE(F326," ","'',''"),"''],''description'':''",G326,"'',''publish'':",H326,"}"),"''","""")'

$r = 327
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 327
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-small-outline-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'people'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A327,",''title'':''",C327,"'',''group'':''",D327,"'',''area'':''",E327,"'',''keywords'':[''",SUBSTITUTE(F327," ","'',''"),"''],''description'':''",G327,"'',''publish'':",H327,"}"),"''","""")'

$r = 328
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 328
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-solid-icons-1-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'icon'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A328,",''title'':''",C328,"'',''group'':''",D328,"'',''area'':''",E328,"'',''keywords'':[''",SUBSTITUTE(F328," ","'',''"),"''],''description'':''",G328,"'',''publish'':",H328,"}"),"''","""")'

$r = 329
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 329
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-solid-icons-2-purple-dark'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'icon'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A329,",''title'':''",C329,"'',''group'':''",D329,"'',''area'':''",E329,"'',''keywords'':[''",SUBSTITUTE(F329," ","'',''"),"''],''description'':''",G329,"'',''publish'':",H329,"}"),"''","""")'

$r = 330
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 330
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-solid-icons-2-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'icon'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A330,",''title'':''",C330,"'',''group'':''",D330,"'',''area'':''",E330,"'',''keywords'':[''",SUBSTITUTE(F330," ","'',''"),"''],''description'':''",G330,"'',''publish'':",H330,"}"),"''","""")'

$r = 331
$ws.Rows.Item($r).RowHeight = 100
$ws.Cells.Item($r, 1).Value = 331
$ws.Cells.Item($r, 1).NumberFormat = "0000"
$ws.Cells.Item($r, 3).Value = 'banner-triangulation-purple'
$ws.Cells.Item($r, 4).Value = 'devblog'
$ws.Cells.Item($r, 5).Value = 'banner'
$ws.Cells.Item($r, 6).Value = 'triangle shape geometry'
$ws.Cells.Item($r, 7).Value = 'Background images used for DevBlog banner area.'
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Formula = '=SUBSTITUTE(_xlfn.CONCAT("{''id'':",A331,",''title'':''",C331,"'',''group'':''",D331,"'',''area'':''",E331,"'',''keywords'':[''",SUBSTITUTE(F331," ","'',''"),"''],''description'':''",G331,"'',''publish'':",H331,"}"),"''","""")'

$ws.Activate()
$ws.Range("I307:I331").Select()
